{"js": "// Splash.docx edit: replace the title paragraph's pict/VML canvas drawing\n// with the updated \"task icon\" artwork (resized canvas, new rectangle\n// backdrop, new \"7\" WordArt badge, repositioned image + \"Sequence\"/\"py\"\n// WordArt) and apply the new Algerian/dark-red run formatting used by the\n// drawing's anchor runs.\n//\n// The legacy VML (v:group/v:shape/v:textpath, i.e. a <w:pict> canvas) is not\n// modelled by the Word JS Shape API (Shape only exposes outer geometry), so\n// the paragraph is replaced in one shot with Paragraph.insertOoxml(...,\n// \"Replace\") carrying the fully-updated paragraph as a FlatOpc OOXML\n// package. This lands the new markup byte-for-byte as authored while every\n// other part of the package (styles, image relationship rId4, etc.) is left\n// untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph that still carries the original drawing (w:pict);\n// robust even if the document gains/loses unrelated paragraphs.\nlet target = null;\nfor (const paragraph of paragraphs.items) {\n  const ooxml = paragraph.getOoxml();\n  await context.sync();\n  if (ooxml.value.includes(\"w:pict\")) {\n    target = paragraph;\n    break;\n  }\n}\nif (!target) {\n  target = paragraphs.items[0];\n}\n\nconst flatOpcPackage = \"<?xml version=\\\"1.0\\\" standalone=\\\"yes\\\"?><?mso-application progid=\\\"Word.Document\\\"?><pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/_rels/.rels\\\" pkg:contentType=\\\"application/vnd.openxmlformats-package.relationships+xml\\\" pkg:padding=\\\"512\\\"><pkg:xmlData><Relationships xmlns=\\\"http://schemas.openxmlformats.org/package/2006/relationships\\\"><Relationship Id=\\\"rId1\\\" Type=\\\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\\\" Target=\\\"word/document.xml\\\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><?xml version=\\\"1.0\\\" encoding=\\\"UTF-8\\\" standalone=\\\"yes\\\"?>\\n<w:document xmlns:ve=\\\"http://schemas.openxmlformats.org/markup-compatibility/2006\\\" xmlns:o=\\\"urn:schemas-microsoft-com:office:office\\\" xmlns:r=\\\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\\\" xmlns:m=\\\"http://schemas.openxmlformats.org/officeDocument/2006/math\\\" xmlns:v=\\\"urn:schemas-microsoft-com:vml\\\" xmlns:wp=\\\"http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing\\\" xmlns:w10=\\\"urn:schemas-microsoft-com:office:word\\\" xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\" xmlns:wne=\\\"http://schemas.microsoft.com/office/word/2006/wordml\\\"><w:body><w:p w:rsidR=\\\"00F622E7\\\" w:rsidRDefault=\\\"00775570\\\"><w:r><w:rPr><w:rFonts w:ascii=\\\"Algerian\\\" w:hAnsi=\\\"Algerian\\\"/><w:color w:val=\\\"C00000\\\"/></w:rPr></w:r><w:r w:rsidR=\\\"006E22AA\\\"><w:rPr><w:rFonts w:ascii=\\\"Algerian\\\" w:hAnsi=\\\"Algerian\\\"/><w:color w:val=\\\"C00000\\\"/></w:rPr><w:pict><v:group id=\\\"_x0000_s1027\\\" editas=\\\"canvas\\\" style=\\\"width:256.55pt;height:180pt;mso-position-horizontal-relative:char;mso-position-vertical-relative:line\\\" coordorigin=\\\"1417,1710\\\" coordsize=\\\"5131,3600\\\"><o:lock v:ext=\\\"edit\\\" aspectratio=\\\"t\\\"/><v:shapetype id=\\\"_x0000_t75\\\" coordsize=\\\"21600,21600\\\" o:spt=\\\"75\\\" o:preferrelative=\\\"t\\\" path=\\\"m@4@5l@4@11@9@11@9@5xe\\\" filled=\\\"f\\\" stroked=\\\"f\\\"><v:stroke joinstyle=\\\"miter\\\"/><v:formulas><v:f eqn=\\\"if lineDrawn pixelLineWidth 0\\\"/><v:f eqn=\\\"sum @0 1 0\\\"/><v:f eqn=\\\"sum 0 0 @1\\\"/><v:f eqn=\\\"prod @2 1 2\\\"/><v:f eqn=\\\"prod @3 21600 pixelWidth\\\"/><v:f eqn=\\\"prod @3 21600 pixelHeight\\\"/><v:f eqn=\\\"sum @0 0 1\\\"/><v:f eqn=\\\"prod @6 1 2\\\"/><v:f eqn=\\\"prod @7 21600 pixelWidth\\\"/><v:f eqn=\\\"sum @8 21600 0\\\"/><v:f eqn=\\\"prod @7 21600 pixelHeight\\\"/><v:f eqn=\\\"sum @10 21600 0\\\"/></v:formulas><v:path o:extrusionok=\\\"f\\\" gradientshapeok=\\\"t\\\" o:connecttype=\\\"rect\\\"/><o:lock v:ext=\\\"edit\\\" aspectratio=\\\"t\\\"/></v:shapetype><v:shape id=\\\"_x0000_s1026\\\" type=\\\"#_x0000_t75\\\" style=\\\"position:absolute;left:1417;top:1710;width:5131;height:3600\\\" o:preferrelative=\\\"f\\\"><v:fill o:detectmouseclick=\\\"t\\\"/><v:path o:extrusionok=\\\"t\\\" o:connecttype=\\\"none\\\"/><o:lock v:ext=\\\"edit\\\" text=\\\"t\\\"/></v:shape><v:rect id=\\\"_x0000_s1031\\\" style=\\\"position:absolute;left:1417;top:1710;width:5131;height:3600\\\" fillcolor=\\\"white [3212]\\\" stroked=\\\"f\\\" strokecolor=\\\"#95b3d7 [1940]\\\" strokeweight=\\\"1pt\\\"><v:fill color2=\\\"#c6d9f1 [671]\\\" focusposition=\\\".5,.5\\\" focussize=\\\"\\\" focus=\\\"50%\\\" type=\\\"gradient\\\"/><v:shadow type=\\\"perspective\\\" color=\\\"#243f60 [1604]\\\" opacity=\\\".5\\\" offset=\\\"1pt\\\" offset2=\\\"-3pt\\\"/></v:rect><v:shapetype id=\\\"_x0000_t136\\\" coordsize=\\\"21600,21600\\\" o:spt=\\\"136\\\" adj=\\\"10800\\\" path=\\\"m@7,l@8,m@5,21600l@6,21600e\\\"><v:formulas><v:f eqn=\\\"sum #0 0 10800\\\"/><v:f eqn=\\\"prod #0 2 1\\\"/><v:f eqn=\\\"sum 21600 0 @1\\\"/><v:f eqn=\\\"sum 0 0 @2\\\"/><v:f eqn=\\\"sum 21600 0 @3\\\"/><v:f eqn=\\\"if @0 @3 0\\\"/><v:f eqn=\\\"if @0 21600 @1\\\"/><v:f eqn=\\\"if @0 0 @2\\\"/><v:f eqn=\\\"if @0 @4 21600\\\"/><v:f eqn=\\\"mid @5 @6\\\"/><v:f eqn=\\\"mid @8 @5\\\"/><v:f eqn=\\\"mid @7 @8\\\"/><v:f eqn=\\\"mid @6 @7\\\"/><v:f eqn=\\\"sum @6 0 @5\\\"/></v:formulas><v:path textpathok=\\\"t\\\" o:connecttype=\\\"custom\\\" o:connectlocs=\\\"@9,0;@10,10800;@11,21600;@12,10800\\\" o:connectangles=\\\"270,180,90,0\\\"/><v:textpath on=\\\"t\\\" fitshape=\\\"t\\\"/><v:handles><v:h position=\\\"#0,bottomRight\\\" xrange=\\\"6629,14971\\\"/></v:handles><o:lock v:ext=\\\"edit\\\" text=\\\"t\\\" shapetype=\\\"t\\\"/></v:shapetype><v:shape id=\\\"_x0000_s1032\\\" type=\\\"#_x0000_t136\\\" style=\\\"position:absolute;left:5134;top:3790;width:690;height:825\\\" fillcolor=\\\"#974706 [1609]\\\" stroked=\\\"f\\\"><v:fill color2=\\\"#aaa\\\"/><v:shadow on=\\\"t\\\" color=\\\"#4d4d4d\\\" opacity=\\\"52429f\\\" offset=\\\",3pt\\\"/><v:textpath style=\\\"font-family:&quot;Arial Black&quot;;v-text-spacing:78650f;v-text-kern:t\\\" trim=\\\"t\\\" fitpath=\\\"t\\\" string=\\\"7\\\"/></v:shape><v:shape id=\\\"_x0000_s1028\\\" type=\\\"#_x0000_t75\\\" style=\\\"position:absolute;left:2908;top:3016;width:2035;height:2036\\\"><v:imagedata r:id=\\\"rId4\\\" o:title=\\\"\\\"/></v:shape><v:shape id=\\\"_x0000_s1030\\\" type=\\\"#_x0000_t136\\\" style=\\\"position:absolute;left:2908;top:2191;width:2731;height:825\\\" fillcolor=\\\"#0070c0\\\" stroked=\\\"f\\\" strokecolor=\\\"green\\\"><v:shadow on=\\\"t\\\" type=\\\"perspective\\\" color=\\\"#c7dfd3\\\" opacity=\\\"52429f\\\" origin=\\\"-.5,-.5\\\" offset=\\\"-3pt,-17pt\\\" offset2=\\\"46pt,38pt\\\" matrix=\\\"1.25,,,1.25\\\"/><v:textpath style=\\\"font-family:&quot;Cambria&quot;;font-weight:bold;v-text-kern:t\\\" trim=\\\"t\\\" fitpath=\\\"t\\\" string=\\\"S\\u00e9quence\\\"/></v:shape><v:shape id=\\\"_x0000_s1029\\\" type=\\\"#_x0000_t136\\\" style=\\\"position:absolute;left:2157;top:2415;width:690;height:825;rotation:-1047489fd\\\" fillcolor=\\\"#4e6128 [1606]\\\" stroked=\\\"f\\\" strokecolor=\\\"green\\\"><v:shadow type=\\\"perspective\\\" color=\\\"#c7dfd3\\\" opacity=\\\"52429f\\\" origin=\\\"-.5,-.5\\\" offset=\\\"-15pt,-12pt\\\" offset2=\\\"22pt,48pt\\\" matrix=\\\"1.25,,,1.25\\\"/><v:textpath style=\\\"font-family:&quot;Times New Roman&quot;;v-text-kern:t\\\" trim=\\\"t\\\" fitpath=\\\"t\\\" string=\\\"py\\\"/></v:shape><w10:wrap type=\\\"none\\\"/><w10:anchorlock/></v:group></w:pict></w:r></w:p><w:sectPr w:rsidR=\\\"00F622E7\\\" w:rsidSect=\\\"00F622E7\\\"><w:pgSz w:w=\\\"11906\\\" w:h=\\\"16838\\\"/><w:pgMar w:top=\\\"1417\\\" w:right=\\\"1417\\\" w:bottom=\\\"1417\\\" w:left=\\\"1417\\\" w:header=\\\"708\\\" w:footer=\\\"708\\\" w:gutter=\\\"0\\\"/><w:cols w:space=\\\"708\\\"/><w:docGrid w:linePitch=\\\"360\\\"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\ntarget.insertOoxml(flatOpcPackage, \"Replace\");\nawait context.sync();\n", "ps1": "# Splash.docx edit: replace the title paragraph's pict/VML canvas drawing\n# with the updated \"task icon\" artwork (resized canvas, new rectangle\n# backdrop, new \"7\" WordArt badge, repositioned image + \"Sequence\"/\"py\"\n# WordArt) and apply the new Algerian/dark-red run formatting used by the\n# drawing's anchor runs. The whole paragraph is swapped in one shot via\n# Range.InsertXML so the legacy VML tree (not modelled by the Shape OM)\n# lands byte-for-byte as authored.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that still carries the original drawing (w:pict);\n# robust even if the document gains/loses unrelated paragraphs.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.WordOpenXML -match \"w:pict\") {\n        $target = $p\n        break\n    }\n}\nif (-not $target) { $target = $d.Paragraphs(1) }\n$r = $target.Range\n\n$xml = @'\n<w:p w:rsidR=\"00F622E7\" w:rsidRDefault=\"00775570\" xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:v=\"urn:schemas-microsoft-com:vml\" xmlns:o=\"urn:schemas-microsoft-com:office:office\" xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\" xmlns:w10=\"urn:schemas-microsoft-com:office:word\"><w:r><w:rPr><w:rFonts w:ascii=\"Algerian\" w:hAnsi=\"Algerian\"/><w:color w:val=\"C00000\"/></w:rPr></w:r><w:r w:rsidR=\"006E22AA\"><w:rPr><w:rFonts w:ascii=\"Algerian\" w:hAnsi=\"Algerian\"/><w:color w:val=\"C00000\"/></w:rPr><w:pict><v:group id=\"_x0000_s1027\" editas=\"canvas\" style=\"width:256.55pt;height:180pt;mso-position-horizontal-relative:char;mso-position-vertical-relative:line\" coordorigin=\"1417,1710\" coordsize=\"5131,3600\"><o:lock v:ext=\"edit\" aspectratio=\"t\"/><v:shapetype id=\"_x0000_t75\" coordsize=\"21600,21600\" o:spt=\"75\" o:preferrelative=\"t\" path=\"m@4@5l@4@11@9@11@9@5xe\" filled=\"f\" stroked=\"f\"><v:stroke joinstyle=\"miter\"/><v:formulas><v:f eqn=\"if lineDrawn pixelLineWidth 0\"/><v:f eqn=\"sum @0 1 0\"/><v:f eqn=\"sum 0 0 @1\"/><v:f eqn=\"prod @2 1 2\"/><v:f eqn=\"prod @3 21600 pixelWidth\"/><v:f eqn=\"prod @3 21600 pixelHeight\"/><v:f eqn=\"sum @0 0 1\"/><v:f eqn=\"prod @6 1 2\"/><v:f eqn=\"prod @7 21600 pixelWidth\"/><v:f eqn=\"sum @8 21600 0\"/><v:f eqn=\"prod @7 21600 pixelHeight\"/><v:f eqn=\"sum @10 21600 0\"/></v:formulas><v:path o:extrusionok=\"f\" gradientshapeok=\"t\" o:connecttype=\"rect\"/><o:lock v:ext=\"edit\" aspectratio=\"t\"/></v:shapetype><v:shape id=\"_x0000_s1026\" type=\"#_x0000_t75\" style=\"position:absolute;left:1417;top:1710;width:5131;height:3600\" o:preferrelative=\"f\"><v:fill o:detectmouseclick=\"t\"/><v:path o:extrusionok=\"t\" o:connecttype=\"none\"/><o:lock v:ext=\"edit\" text=\"t\"/></v:shape><v:rect id=\"_x0000_s1031\" style=\"position:absolute;left:1417;top:1710;width:5131;height:3600\" fillcolor=\"white [3212]\" stroked=\"f\" strokecolor=\"#95b3d7 [1940]\" strokeweight=\"1pt\"><v:fill color2=\"#c6d9f1 [671]\" focusposition=\".5,.5\" focussize=\"\" focus=\"50%\" type=\"gradient\"/><v:shadow type=\"perspective\" color=\"#243f60 [1604]\" opacity=\".5\" offset=\"1pt\" offset2=\"-3pt\"/></v:rect><v:shapetype id=\"_x0000_t136\" coordsize=\"21600,21600\" o:spt=\"136\" adj=\"10800\" path=\"m@7,l@8,m@5,21600l@6,21600e\"><v:formulas><v:f eqn=\"sum #0 0 10800\"/><v:f eqn=\"prod #0 2 1\"/><v:f eqn=\"sum 21600 0 @1\"/><v:f eqn=\"sum 0 0 @2\"/><v:f eqn=\"sum 21600 0 @3\"/><v:f eqn=\"if @0 @3 0\"/><v:f eqn=\"if @0 21600 @1\"/><v:f eqn=\"if @0 0 @2\"/><v:f eqn=\"if @0 @4 21600\"/><v:f eqn=\"mid @5 @6\"/><v:f eqn=\"mid @8 @5\"/><v:f eqn=\"mid @7 @8\"/><v:f eqn=\"mid @6 @7\"/><v:f eqn=\"sum @6 0 @5\"/></v:formulas><v:path textpathok=\"t\" o:connecttype=\"custom\" o:connectlocs=\"@9,0;@10,10800;@11,21600;@12,10800\" o:connectangles=\"270,180,90,0\"/><v:textpath on=\"t\" fitshape=\"t\"/><v:handles><v:h position=\"#0,bottomRight\" xrange=\"6629,14971\"/></v:handles><o:lock v:ext=\"edit\" text=\"t\" shapetype=\"t\"/></v:shapetype><v:shape id=\"_x0000_s1032\" type=\"#_x0000_t136\" style=\"position:absolute;left:5134;top:3790;width:690;height:825\" fillcolor=\"#974706 [1609]\" stroked=\"f\"><v:fill color2=\"#aaa\"/><v:shadow on=\"t\" color=\"#4d4d4d\" opacity=\"52429f\" offset=\",3pt\"/><v:textpath style=\"font-family:&quot;Arial Black&quot;;v-text-spacing:78650f;v-text-kern:t\" trim=\"t\" fitpath=\"t\" string=\"7\"/></v:shape><v:shape id=\"_x0000_s1028\" type=\"#_x0000_t75\" style=\"position:absolute;left:2908;top:3016;width:2035;height:2036\"><v:imagedata r:id=\"rId4\" o:title=\"\"/></v:shape><v:shape id=\"_x0000_s1030\" type=\"#_x0000_t136\" style=\"position:absolute;left:2908;top:2191;width:2731;height:825\" fillcolor=\"#0070c0\" stroked=\"f\" strokecolor=\"green\"><v:shadow on=\"t\" type=\"perspective\" color=\"#c7dfd3\" opacity=\"52429f\" origin=\"-.5,-.5\" offset=\"-3pt,-17pt\" offset2=\"46pt,38pt\" matrix=\"1.25,,,1.25\"/><v:textpath style=\"font-family:&quot;Cambria&quot;;font-weight:bold;v-text-kern:t\" trim=\"t\" fitpath=\"t\" string=\"S\u00e9quence\"/></v:shape><v:shape id=\"_x0000_s1029\" type=\"#_x0000_t136\" style=\"position:absolute;left:2157;top:2415;width:690;height:825;rotation:-1047489fd\" fillcolor=\"#4e6128 [1606]\" stroked=\"f\" strokecolor=\"green\"><v:shadow type=\"perspective\" color=\"#c7dfd3\" opacity=\"52429f\" origin=\"-.5,-.5\" offset=\"-15pt,-12pt\" offset2=\"22pt,48pt\" matrix=\"1.25,,,1.25\"/><v:textpath style=\"font-family:&quot;Times New Roman&quot;;v-text-kern:t\" trim=\"t\" fitpath=\"t\" string=\"py\"/></v:shape><w10:wrap type=\"none\"/><w10:anchorlock/></v:group></w:pict></w:r></w:p>\n'@\n\n$r.InsertXML($xml)\n"}
